$p = $ppt.ActivePresentation

# --- Remove the second slide entirely ---
$p.Slides.Item(2).Delete()

$s1 = $p.Slides.Item(1)

# --- Remove the two picture shapes from slide 1, keep only the text box ---
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Type -eq 13) {
        $shp.Delete()
    }
}

# --- The remaining shape is the text box; rename, reposition & resize it ---
$tb = $s1.Shapes.Item(1)
$tb.Name = "TextBox 1"
$tb.Left = 3419061 / 12700
$tb.Top = 864704 / 12700
$tb.Width = 1421296 / 12700
$tb.Height = 923330 / 12700

# --- Replace the text content with the new paragraphs ---
$tr = $tb.TextFrame.TextRange
$tr.Text = "# palabras`r `r" + [char]0x2018
$tr2 = $tr.InsertAfter("asd")
$tr3 = $tr2.InsertAfter([char]0x2019)
